# Generate Report for Handback
# Marks the two handed-back files (218549ab..., c35c17c7...) as synced and
# fills in their Latest Target File / Latest Handback File / Latest Handback
# DateTime columns on the zh-cn and de-de locale sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/72e703a3c437dc03bf5a842ea2d5d5b6be83deec/e2e/218549ab-b967-417e-b842-e836e0710af4.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/72e703a3c437dc03bf5a842ea2d5d5b6be83deec/e2e/c35c17c7-d1d0-4b78-bd4a-216cfb0dbcd1.md"

$mdName1 = "218549ab-b967-417e-b842-e836e0710af4.md"
$mdName2 = "c35c17c7-d1d0-4b78-bd4a-216cfb0dbcd1.md"

# ----------------------------------------------------------------------
# Overview sheet: mark both locales as handed back / in sync for both files
# ----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# ----------------------------------------------------------------------
# zh-cn sheet: mark status as handed back + fill in target/handback file
# and handback datetime
# ----------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("J2").Value = "218549ab-b967-417e-b842-e836e0710af4.4dc40e804199b899e728ad83da9ee7b6129ffad1.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-10-18 05:10:23"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl1, "", "", $mdName1)

$zhcn.Range("J3").Value = "c35c17c7-d1d0-4b78-bd4a-216cfb0dbcd1.490c629aa192833ad0e1dd6b486e929e57bd1d11.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-10-18 05:10:23"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl2, "", "", $mdName2)

# ----------------------------------------------------------------------
# de-de sheet: mark status as handed back + fill in target/handback file
# and handback datetime
# ----------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("J2").Value = "218549ab-b967-417e-b842-e836e0710af4.4dc40e804199b899e728ad83da9ee7b6129ffad1.de-de.xlf"
$dede.Range("K2").Value = "2016-10-18 05:10:47"
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl1, "", "", $mdName1)

$dede.Range("J3").Value = "c35c17c7-d1d0-4b78-bd4a-216cfb0dbcd1.490c629aa192833ad0e1dd6b486e929e57bd1d11.de-de.xlf"
$dede.Range("K3").Value = "2016-10-18 05:10:47"
$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl2, "", "", $mdName2)
